# Applies the etable_4.docx edit:
#  - heading text: drop the "Produce a table from " prefix
#  - translate several table labels/values to Chinese
#  - bump two numeric estimates in the "Turn Circle" row
#  - drop the stray Normal->DocDefaults basedOn link and the now-unused
#    DocDefaults style from styles.xml

$d = $word.ActiveDocument

# --- Heading paragraph -----------------------------------------------
# Rebuild paragraph 1 via InsertParagraphBefore/Delete so the fresh
# paragraph doesn't carry forward the legacy w14:paraId/w14:textId pair.
$p1 = $d.Paragraphs(1)
$r = $p1.Range
$r.Collapse(1)
$r.InsertParagraphBefore()
$newHeading = $d.Paragraphs(1)
$oldHeading = $d.Paragraphs(2)
$newHeading.Range.Text = "-estimates table-"
$oldHeading.Range.Delete()

# --- Table label / value translations ---------------------------------
$d.Content.Find.Execute("model1", $true, $false, $false, $false, $false, `
    $true, 1, $false, "模型1", 2) | Out-Null
$d.Content.Find.Execute("model2", $true, $false, $false, $false, $false, `
    $true, 1, $false, "模型2", 2) | Out-Null
$d.Content.Find.Execute("Weight (lbs.)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "重量(公斤)", 2) | Out-Null
$d.Content.Find.Execute("0.0014", $true, $false, $false, $false, $false, `
    $true, 1, $false, "0.0030", 2) | Out-Null
$d.Content.Find.Execute("0.0013", $true, $false, $false, $false, $false, `
    $true, 1, $false, "0.0028", 2) | Out-Null
$d.Content.Find.Execute("Gear Ratio", $true, $false, $false, $false, $false, `
    $true, 1, $false, "变速比", 2) | Out-Null
$d.Content.Find.Execute("Turn Circle (ft.) ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "转弯半径(米) ", 2) | Out-Null
$d.Content.Find.Execute("0.0243", $true, $false, $false, $false, $false, `
    $true, 1, $false, "0.0798", 2) | Out-Null
$d.Content.Find.Execute("0.0613", $true, $false, $false, $false, $false, `
    $true, 1, $false, "0.2010", 2) | Out-Null
$d.Content.Find.Execute("Car type", $true, $false, $false, $false, $false, `
    $true, 1, $false, "国籍", 2) | Out-Null

# --- styles.xml cleanup -------------------------------------------------
# Normal used to be based on a "DocDefaults" style; unlink it, then the
# (now unused) DocDefaults style can be deleted outright.
$normalStyle = $d.Styles.Item("Normal")
$normalStyle.BaseStyle = ""
$docDefaultsStyle = $d.Styles.Item("DocDefaults")
$docDefaultsStyle.Delete()
